# Insert a new weekly price record as row 153 on the (single) worksheet,
# pushing the existing rows 153-217 down to 154-218 (dimension grows from
# A1:R217 to A1:R218).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(153).Insert()

$ws.Cells.Item(153, 1).Value  = 5
$ws.Cells.Item(153, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(153, 3).Value  = "Maule"
$ws.Cells.Item(153, 4).Value  = 44636
$ws.Cells.Item(153, 5).Value  = 7
$ws.Cells.Item(153, 6).Value  = 100112021
$ws.Cells.Item(153, 7).Value  = "Ají"
$ws.Cells.Item(153, 8).Value  = "Cristal"
$ws.Cells.Item(153, 9).Value  = "Primera"
$ws.Cells.Item(153, 10).Value = 100
$ws.Cells.Item(153, 11).Value = 15000
$ws.Cells.Item(153, 12).Value = 15000
$ws.Cells.Item(153, 13).Value = 15000
$ws.Cells.Item(153, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 600
$ws.Cells.Item(153, 17).Value = 25
$ws.Cells.Item(153, 18).Value = "Hortaliza"
